# #2 update mock up file compare
#
# Renames the sheet "SampleTest" -> "SampleTest1" and moves the active
# selection from the whole-sheet range (A1:XFD1048576) to D25, mirroring
# the <selection activeCell="D25" sqref="D25"/> seen in the target
# worksheet XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "SampleTest1"

# Move/collapse the selection onto D25.
$ws.Range("D25").Select()
